$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.139.43'
$ws.Range('E2').Value = '  +2.59%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.403.43'
$ws.Range('E3').Value = '  +2.83%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.84'
$ws.Range('E5').Value = '  +2.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.99'
$ws.Range('E6').Value = '  +4.54%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +0.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.401.55'
$ws.Range('E9').Value = '  +2.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.104'
$ws.Range('E10').Value = '  +2.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.67'
$ws.Range('E11').Value = '  +2.67%  '
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.349'
$ws.Range('E13').Value = '  +3.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.83'
$ws.Range('E14').Value = '  +8.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.829.94'
$ws.Range('E15').Value = '  +2.77%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '62.060.39'
$ws.Range('E16').Value = '  +2.58%  '
$ws.Range('E17').Value = '  +4.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.386.79'
$ws.Range('E18').Value = '  +2.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.04'
$ws.Range('E19').Value = '  +3.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '344.01'
$ws.Range('E20').Value = '  +8.93%  '
$ws.Range('E21').Value = '  +1.43%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.89'
$ws.Range('E22').Value = '  +2.95%  '
$ws.Range('E23').Value = '  +0.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.21'
$ws.Range('E24').Value = '  +1.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.173'
$ws.Range('E25').Value = '  +1.40%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.38'
$ws.Range('E27').Value = '  +6.33%  '
$ws.Range('E28').Value = '  +10.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.39'
$ws.Range('E29').Value = '  +14.72%  '
$ws.Range('E30').Value = '  +3.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0769'
$ws.Range('E31').Value = '  +3.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.36'
$ws.Range('E32').Value = '  +6.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '171.33'
$ws.Range('E33').Value = '  -1.40%  '
$ws.Range('E34').Value = '  +2.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.394'
$ws.Range('E35').Value = '  +3.23%  '
$ws.Range('E36').Value = '  +3.17%  '
$ws.Range('E37').Value = '  +10.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.998'
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '356.99'
$ws.Range('E39').Value = '  +9.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.67'
$ws.Range('E41').Value = '  +8.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.10'
$ws.Range('E42').Value = '  +2.46%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '144.27'
$ws.Range('E43').Value = '  +3.03%  '
$ws.Range('E44').Value = '  +5.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.65'
$ws.Range('E45').Value = '  +6.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0964'
$ws.Range('E46').Value = '  +2.00%  '
$ws.Range('E47').Value = '  +4.18%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.584'
$ws.Range('E48').Value = '  +3.79%  '
$ws.Range('E49').Value = '  +3.54%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.84'
$ws.Range('E50').Value = '  +4.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0215'
$ws.Range('E51').Value = '  -4.35%  '
